$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$pi = 3.141592653589793

for ($r = 1; $r -le 17; $r++) {
    $b = $ws.Cells.Item($r, 2).Value2
    $c = $ws.Cells.Item($r, 3).Value2
    $d = $ws.Cells.Item($r, 4).Value2
    $e = $ws.Cells.Item($r, 5).Value2
    $f = $ws.Cells.Item($r, 6).Value2

    $ws.Cells.Item($r, 2).Value = -$b
    $ws.Cells.Item($r, 3).Value = ($pi / 2) - $c
    $ws.Cells.Item($r, 4).Value = -$d
    $ws.Cells.Item($r, 5).Value = ($pi / 2) - $e
    $ws.Cells.Item($r, 6).Value = -$f
}
